$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-reported daily error counts for 1/6/2026 and 1/7/2026
# (row 68 -> serial 46028, row 69 -> serial 46029), matching the existing
# Date / Total Count / Session Timeout Errors / Errors Requiring Analysis
# columns already used by the rest of the sheet.
$ws.Range("A68").Value = 46028
$ws.Range("B68").Value = 618
$ws.Range("C68").Value = 590
$ws.Range("D68").Value = 28

$ws.Range("A69").Value = 46029
$ws.Range("B69").Value = 554
$ws.Range("C69").Value = 527
$ws.Range("D69").Value = 27

# Update the view so the window is scrolled down to show the new rows and
# the active cell/selection matches what was saved (G67).
$excel.ActiveWindow.TopLeftCell = $ws.Range("A56")
$ws.Range("G67").Select() | Out-Null
